$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row 1: build "Iteration_1" (E1:G1) and "Iteration_2" (H1:J1) blocks,
#     mirroring the existing "Standalone" (B1:D1) block's styling/layout. ---

# Seed E1/H1 (label cell) and F1,G1 / I1,J1 (blank-but-styled cells) from the
# existing A1 (styled label) and C1/D1 (styled, empty) cells so the style (s="1")
# carries over exactly as it does for B1:D1.
$ws.Range("A1").Copy($ws.Range("E1"))
$ws.Range("C1").Copy($ws.Range("F1"))
$ws.Range("D1").Copy($ws.Range("G1"))

$ws.Range("A1").Copy($ws.Range("H1"))
$ws.Range("C1").Copy($ws.Range("I1"))
$ws.Range("D1").Copy($ws.Range("J1"))

$ws.Range("E1").Value = "Iteration_1"
$ws.Range("H1").Value = "Iteration_2"

# Merge the new header blocks the same way B1:D1 is merged, then reassert the
# thin-box/bold/centered formatting (merging nudges the per-cell border
# variant) so every cell in the block resolves back to the original shared
# style used across row 1.
$ws.Range("E1:G1").Merge()
$ws.Range("E1:G1").Borders.LineStyle = 1

$ws.Range("H1:J1").Merge()
$ws.Range("H1:J1").Borders.LineStyle = 1

# --- Row 2 (Interval sub-header: 2030 / 2040 / 2050) ---
# Copy straight from B2:D2 so text-typed values + style (s="1") match exactly.
$ws.Range("B2:D2").Copy($ws.Range("E2"))
$ws.Range("B2:D2").Copy($ws.Range("H2"))

# --- Data rows 4-16, columns E:J ---

$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0

$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0

$ws.Range("B6").Value = -0.000000002235549776991508
$ws.Range("E6").Value = 1401265.181013603
$ws.Range("F6").Value = 602827.8937104597
$ws.Range("G6").Value = 221185.4197186728
$ws.Range("H6").Value = 1401265.18101358
$ws.Range("I6").Value = 603806.4382143666
$ws.Range("J6").Value = 222160.8795700106

$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0

$ws.Range("C8").Value = 1205602.936901787
$ws.Range("D8").Value = 1205602.936901779
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 832311.6761443499
$ws.Range("G8").Value = 894271.2363347108
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 831090.3558850527
$ws.Range("J8").Value = 892473.7149865876

$ws.Range("B9").Value = 1717310
$ws.Range("C9").Value = 15814.1753639676
$ws.Range("D9").Value = -0.0000000001136217926629919
$ws.Range("E9").Value = 316335.4517647082
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = -0.00000000006340789397460244
$ws.Range("H9").Value = 316335.4517647056
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = -0.00000000007024058041087535

$ws.Range("C10").Value = 495892.8877342517
$ws.Range("D10").Value = 511707.0630982184
$ws.Range("E10").Value = 0
$ws.Range("F10").Value = 282295.4611157335
$ws.Range("G10").Value = 601899.2194410748
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 282538.4398285086
$ws.Range("J10").Value = 602721.4832554583

$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 0
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0

$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 0

$ws.Range("E13").Value = 0
$ws.Range("F13").Value = 0
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 0

$ws.Range("E14").Value = 0
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 0

$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 0

$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 0
$ws.Range("G16").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 0
$ws.Range("J16").Value = 0
